$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the commit hash for this release (as a formula yielding a text value)
$ws.Range("B4").Formula = '="514716e1"'
$ws.Range("B4").NumberFormat = "0.00E+00"
$ws.Range("B4").HorizontalAlignment = -4108

# New checklist row 14: NUGET version marked before commit?
$ws.Range("B14:J14").HorizontalAlignment = -4108
$ws.Range("A14").Value = "No"
$ws.Range("B14:J14").Merge()
$ws.Range("B14").Value = "Was version properly marked in NUGET configuration before commit?"

# New checklist row 15: project/NUGET configuration version updated?
$ws.Range("B15:J15").HorizontalAlignment = -4108
$ws.Range("A15").Value = "Yes"
$ws.Range("B15:J15").Merge()
$ws.Range("B15").Value = "Has version in project configuration and NUGET config been updated?"

[void]$ws.Range("B14:J14").Select()
